# Auto-generated Excel COM-interop script that refreshes the crypto price /
# volume table (columns D and E) for rows 2-51, and swaps the FraxShare /
# Cronos rows (45-46, including their Coin name and Link columns) to match
# the latest data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new Price values (column D) are plain decimal numbers (e.g.
# "2.06", "0.0920", "0.0000101"). If assigned as-is, Excel auto-detects them
# as numbers and silently reformats/loses the exact text (e.g. trailing
# zeros). The source sheet stores these as literal text, so we temporarily
# force each of those cells to the "@" (Text) number format before writing
# the value, then restore the default formatting afterwards so the saved
# styles match the original file.
$numericLookingCells = @(
    "D5"
    "D6"
    "D7"
    "D9"
    "D10"
    "D11"
    "D12"
    "D14"
    "D16"
    "D18"
    "D20"
    "D21"
    "D22"
    "D23"
    "D24"
    "D26"
    "D27"
    "D28"
    "D31"
    "D32"
    "D33"
    "D34"
    "D35"
    "D38"
    "D39"
    "D40"
    "D41"
    "D42"
    "D43"
    "D44"
    "D45"
    "D46"
    "D47"
    "D48"
    "D50"
)
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# --- Update Price (D) / Volume(1h) (E) cells, and the swapped Coin (B) /
# --- Link (C) cells for rows 45-46. ---
$ws.Range("D2").Value = "41.334.27"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "2.177.17"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "238.03"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").Value = "0.609"
$ws.Range("E6").Value = "  -2.39%  "
$ws.Range("D7").Value = "69.82"
$ws.Range("E7").Value = "  -4.34%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.574"
$ws.Range("E9").Value = "  -4.85%  "
$ws.Range("D10").Value = "39.37"
$ws.Range("E10").Value = "  -7.86%  "
$ws.Range("D11").Value = "0.0920"
$ws.Range("E11").Value = "  -3.18%  "
$ws.Range("D12").Value = "54.37"
$ws.Range("E12").Value = "  -5.35%  "
$ws.Range("E13").Value = "  -2.10%  "
$ws.Range("D14").Value = "6.72"
$ws.Range("E14").Value = "  -4.79%  "
$ws.Range("D15").Value = "2.493.99"
$ws.Range("E15").Value = "  -1.95%  "
$ws.Range("D16").Value = "14.09"
$ws.Range("E16").Value = "  -0.72%  "
$ws.Range("D17").Value = "2.162.04"
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("D18").Value = "0.793"
$ws.Range("E18").Value = "  -4.83%  "
$ws.Range("D19").Value = "41.079.56"
$ws.Range("E19").Value = "  -1.66%  "
$ws.Range("D20").Value = "0.0000101"
$ws.Range("E20").Value = "  -6.76%  "
$ws.Range("D21").Value = "70.52"
$ws.Range("E21").Value = "  -3.16%  "
$ws.Range("D22").Value = "5.87"
$ws.Range("E22").Value = "  -4.27%  "
$ws.Range("D23").Value = "225.29"
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("D24").Value = "9.42"
$ws.Range("E24").Value = "  -8.31%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "1.90"
$ws.Range("E26").Value = "  -8.89%  "
$ws.Range("D27").Value = "10.73"
$ws.Range("E27").Value = "  -7.91%  "
$ws.Range("D28").Value = "3.46"
$ws.Range("E28").Value = "  -3.69%  "
$ws.Range("E29").Value = "  -2.80%  "
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").Value = "167.47"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").Value = "19.91"
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("D33").Value = "30.08"
$ws.Range("E33").Value = "  +3.29%  "
$ws.Range("D34").Value = "0.0760"
$ws.Range("E34").Value = "  -4.06%  "
$ws.Range("D35").Value = "5.10"
$ws.Range("E35").Value = "  -9.23%  "
$ws.Range("E36").Value = "  -3.32%  "
$ws.Range("E37").Value = "  -8.56%  "
$ws.Range("D38").Value = "4.06"
$ws.Range("D39").Value = "0.0282"
$ws.Range("E39").Value = "  -5.90%  "
$ws.Range("D40").Value = "2.06"
$ws.Range("E40").Value = "  -2.08%  "
$ws.Range("D41").Value = "11.75"
$ws.Range("E41").Value = "  -13.24%  "
$ws.Range("D42").Value = "5.36"
$ws.Range("E42").Value = "  -4.24%  "
$ws.Range("D43").Value = "58.88"
$ws.Range("E43").Value = "  -9.96%  "
$ws.Range("D44").Value = "0.189"
$ws.Range("E44").Value = "  -4.23%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "0.0970"
$ws.Range("E45").Value = "  -3.12%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "8.23"
$ws.Range("E46").Value = "  -4.71%  "
$ws.Range("D47").Value = "97.33"
$ws.Range("E47").Value = "  -6.01%  "
$ws.Range("D48").Value = "1.08"
$ws.Range("E48").Value = "  -3.18%  "
$ws.Range("E49").Value = "  -3.45%  "
$ws.Range("D50").Value = "2.19"
$ws.Range("E50").Value = "  -7.80%  "
$ws.Range("E51").Value = "  -2.54%  "

# Restore the cells we temporarily switched to Text format back to the
# workbook default (no explicit number format / style).
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "General"
    $ws.Range($cellRef).Style = "Normal"
}

Write-Output "Updated cryptos list: refreshed prices/volumes and reordered FraxShare/Cronos."
